# #39: Add APIs to remove paragraph portion
#
# 1) Refresh the cached "datetimeFigureOut" field text (slide master +
#    every slide layout) from 12/23/2020 -> 12/25/2020.
# 2) On slide 2, shape "TextBox 4", split the first paragraph's single
#    run ("Test-p6") into two portions so a leading bold "id5" portion
#    can be removed later: "id5" (bold) + "-Test-p6" (regular). Widen
#    the (spAutoFit / wrap=none) text box so the new text still fits.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = "12/25/2020"
        }
    }
}

# Slide master's own date placeholder.
Update-DatePlaceholder($p.SlideMaster.Shapes)

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder($layouts.Item($j).Shapes)
}

# Slide 2 / "TextBox 4": widen the box, then split "Test-p6" into a bold
# "id5" portion followed by "-Test-p6".
$slide2 = $p.Slides.Item(2)
$textBox4 = $slide2.Shapes.Item(4)

$textBox4.Width = 124.30216

$tr = $textBox4.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.InsertBefore("id5") | Out-Null

$boldPortion = $textBox4.TextFrame.TextRange.Characters(1, 3)
$boldPortion.Font.Bold = 1

$remainder = $textBox4.TextFrame.TextRange.Characters(4, 7)
$remainder.Text = "-Test-p6"
